$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'sliding knee sleeve'
$ws.Range("A2").Value = 'volleyball spandex shorts'
$ws.Range("A3").Value = 'compression spandex men'
$ws.Range("A4").Value = 'volleyball knee pad'
$ws.Range("A5").Value = 'knee pads for volleyball'
$ws.Range("A6").Value = 'skateboard knee pads'
$ws.Range("A7").Value = 'starter youth compression pants'
$ws.Range("A8").Value = 'mens compression pants marvel'
$ws.Range("A9").Value = 'mouthguard basketball youth'
$ws.Range("A10").Value = 'nike compression leggings'
$ws.Range("A11").Value = 'knee pad sleeves'
$ws.Range("A12").Value = 'jordan flight mens basketball pants'
$ws.Range("A13").Value = 'athletic compression pants'
$ws.Range("A14").Value = 'nike pro compression leggings men'
$ws.Range("A15").Value = 'mens compression tights nike'
$ws.Range("A16").Value = 'ladies compression pants'
$ws.Range("A17").Value = 'mcdavid compression pants'
$ws.Range("A18").Value = 'basketball youth jersey'
$ws.Range("A19").Value = 'protec knee pads'
$ws.Range("A20").Value = 'skate knee pads'
$ws.Range("A21").Value = 'youth knee and elbow pads'
$ws.Range("A22").Value = 'under armour compression tights men'
$ws.Range("A23").Value = 'elbow and knee pads'
$ws.Range("A24").Value = 'adidas tights men'
$ws.Range("A25").Value = 'compression pants tesla'
$ws.Range("A26").Value = 'poc knee pads'
$ws.Range("A27").Value = 'basket knee pads'
$ws.Range("A28").Value = 'exercise knee pad'
$ws.Range("A29").Value = 'mens compression pants adidas'
$ws.Range("A30").Value = 'the rock mens basketball'
$ws.Range("A31").Value = 'men''s tights leggings'
$ws.Range("A32").Value = 'super compression leggings'
$ws.Range("A33").Value = 'under armour compression pants youth boys'
$ws.Range("A34").Value = 'knee pad compression pants'
$ws.Range("A35").Value = 'pantalones con rodilleras para hombre'
$ws.Range("A36").Value = 'basketball tights with knee pads for men'
$ws.Range("A37").Value = 'compression pants men basketball'
$ws.Range("A38").Value = 'knee leggings for basketball'
$ws.Range("A39").Value = 'pants with knee pads for men'
$ws.Range("A40").Value = 'knee pad pants for men'
$ws.Range("A41").Value = 'mens basketball compression pants'
$ws.Range("A42").Value = 'basketball yoga pants'
$ws.Range("A43").Value = 'basketball tights with pads'
$ws.Range("A44").Value = 'sliding pants mens'
$ws.Range("A45").Value = 'basketball compression pants knee pads'
$ws.Range("A46").Value = 'compression with knee pads'
$ws.Range("A47").Value = 'compression knee pads for men'
$ws.Range("A48").Value = 'men basketball compression knee pads'
$ws.Range("A49").Value = 'basketball leggings with pads'
$ws.Range("A50").Value = 'capri pads'
$ws.Range("A51").Value = 'basketball knee pad pants'
$ws.Range("A52").Value = 'tights with pads for men'
$ws.Range("A53").Value = 'mizuno slider knee pad'
$ws.Range("A54").Value = 'mizuno adult slider kneepad'
$ws.Range("A55").Value = 'compression knee leggings'
$ws.Range("A56").Value = '5 pad compression shorts'
$ws.Range("A57").Value = 'mizuno slider kneepad'
$ws.Range("A58").Value = 'youth basketball tights with knee pads'
$ws.Range("A59").Value = 'basketball knee pads for men'
$ws.Range("A60").Value = 'legging pads'
$ws.Range("A61").Value = 'compression capri men pack'
$ws.Range("A62").Value = 'softball knee sliding pad'
$ws.Range("A63").Value = 'basketball tights knee'
$ws.Range("A64").Value = 'compression pads basketball'
$ws.Range("A65").Value = 'leggings with knee pads'
$ws.Range("A66").Value = 'knee pads baseball'
$ws.Range("A67").Value = 'knee length tights men'
$ws.Range("A68").Value = 'basketball tight'
$ws.Range("A69").Value = 'baseball sliding knee pad'
$ws.Range("A70").Value = 'knee compression pants'
$ws.Range("A71").Value = 'baseball sliding pad'
$ws.Range("A72").Value = 'padded compression pants basketball'
$ws.Range("A73").Value = 'compression pants capri men'
$ws.Range("A74").Value = 'basketball men leggings'
$ws.Range("A75").Value = 'knee sliders softball'
$ws.Range("A76").Value = 'mens padded basketball tights'
$ws.Range("A77").Value = 'compression pants mens basketball'
$ws.Range("A78").Value = 'softball knee slider'
$ws.Range("A79").Value = 'mens weightlifting tights'
$ws.Range("A80").Value = 'mcdavid 6446 hex knee pads compression leg sleeve'
$ws.Range("A81").Value = 'softball sliding pad'
$ws.Range("A82").Value = 'pants pad'
$ws.Range("A83").Value = 'knee pads mens basketball'
$ws.Range("A84").Value = 'padded pants men'
$ws.Range("A85").Value = 'basketball legging'
$ws.Range("A86").Value = 'padded compression tights basketball'
$ws.Range("A87").Value = 'youth tights with knee pads'
$ws.Range("A88").Value = 'youth padded tights'
$ws.Range("A89").Value = 'baseball knee pads adult'
$ws.Range("A90").Value = 'mens capri leggings for sports'
$ws.Range("A91").Value = 'basketball tights with knee pads youth boys'
$ws.Range("A92").Value = 'youth compression pants with knee pads'
$ws.Range("A93").Value = 'mens capris pants'
$ws.Range("A94").Value = 'basketball knee pad leggings'
$ws.Range("A95").Value = 'padded pants for basketball'
$ws.Range("A96").Value = 'compression capri pants men'
$ws.Range("A97").Value = 'sliding knee pads'
$ws.Range("A98").Value = 'mens capri compression tights'
$ws.Range("A99").Value = 'athletic below knee'
$ws.Range("A100").Value = 'softball sliding pads'
